# Added tab switching for GHG molecule charts
# - Adds a new "CO2" worksheet after Sheet1 containing the year / ppm
#   columns (copied from Sheet1's year + co2 columns), with matching
#   header/number formatting, and makes it the active tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add the new sheet right after Sheet1 (becomes the last sheet / active tab).
$co2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$co2.Name = "CO2"

# Header row.
$co2.Range("A1").Value = "YEAR"
$co2.Range("B1").Value = "PPM"

# Year / co2-ppm data, mirroring Sheet1!A6:B35 (1984-2013).
$years = 1984,1985,1986,1987,1988,1989,1990,1991,1992,1993,1994,1995,1996,1997,1998,1999,2000,2001,2002,2003,2004,2005,2006,2007,2008,2009,2010,2011,2012,2013
$ppm   = 344.3,345.7,347.2,349,351.3,353,354.1,355.4,356.1,357,358.5,360.3,362,363.3,366,368,369.4,370.9,372.9,375.3,377.1,379.2,381.3,383.1,385.2,386.8,388.9,390.9,393.1,396

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $co2.Cells.Item($row, 1).Value = $years[$i]
    $co2.Cells.Item($row, 2).Value = $ppm[$i]
}

# Match Sheet1's header/number formatting (reuse existing styles).
$sheet1.Range("A1").Copy()
$co2.Range("A1").PasteSpecial(-4122)
$sheet1.Range("B1").Copy()
$co2.Range("B1").PasteSpecial(-4122)

$sheet1.Range("A2").Copy()
$co2.Range("A2:A31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Make CO2 the active tab/selection, matching the tab-switch the commit adds.
[void]$co2.Range("B4").Select()
